# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" table (rows 16-124) is re-keyed in ascending
# chronological order (1607 .. 2507) instead of the previous descending
# order, and a brand new period (2508) is appended as a new last data row.
# The fixed "Valor Mora" amount per period also follows the period
# (<=1808 -> 24640, >=1809 -> 31249), which is why many existing rows'
# F-column value flips even though the row number itself does not move.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Insert a new row at 125. This pushes the blank rows 125-128 down to
#    126-129, and the two signature-block rows (129 -> 130, 130 -> 131)
#    down by one, exactly like the target file (mergeCells follow suit
#    automatically).
# ---------------------------------------------------------------------
$ws.Rows.Item(125).Insert()

# ---------------------------------------------------------------------
# 2. Row 124 used to be the special "last row" (bold border) style; now
#    that row 125 is the new last row, swap the two rows' visual style:
#      - give row 125 the old "last row" look (copied from row 124)
#      - give row 124 the normal data-row look (copied from row 123)
# ---------------------------------------------------------------------
$ws.Range("B124:J124").Copy() | Out-Null
$ws.Range("B125:J125").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B123:J123").Copy() | Out-Null
$ws.Range("B124:J124").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. Build the full ascending list of periods: 1607-1612, 1701-2412,
#    2501-2508 (110 periods total, one more than the 109 before the
#    edit).
# ---------------------------------------------------------------------
$periods = New-Object System.Collections.ArrayList
for ($m = 7; $m -le 12; $m++) {
    [void]$periods.Add("16{0:D2}" -f $m)
}
for ($y = 17; $y -le 24; $y++) {
    for ($m = 1; $m -le 12; $m++) {
        [void]$periods.Add("{0}{1:D2}" -f $y, $m)
    }
}
for ($m = 1; $m -le 8; $m++) {
    [void]$periods.Add("25{0:D2}" -f $m)
}

# ---------------------------------------------------------------------
# 4. Re-key rows 16-125 (the 110 period/data rows) with the new
#    ascending period list and the period-dependent "Valor Mora" amount.
# ---------------------------------------------------------------------
for ($i = 0; $i -lt $periods.Count; $i++) {
    $r = 16 + $i
    $period = [string]$periods[$i]
    if ([int]$period -le 1808) {
        $valor = 24640
    } else {
        $valor = 31249
    }
    $ws.Cells.Item($r, 2).Value2 = "CC"                              # column B - Tipo Doc Trabajador
    $ws.Cells.Item($r, 3).Value2 = "1003197017"                      # column C - N Doc Trabajador
    $ws.Cells.Item($r, 4).Value2 = "GILBERTO LUIS BORNACHERA YEPEZ"  # column D - Nombre Trabajador
    $ws.Cells.Item($r, 5).Value2 = $period      # column E - Periodo Mora
    $ws.Cells.Item($r, 6).Value2 = $valor       # column F - Valor Mora
    $ws.Cells.Item($r, 7).Value2 = 781242       # column G - Salario Basico
}

# ---------------------------------------------------------------------
# 5. Fix up the summary cells driven by the table above.
# ---------------------------------------------------------------------
$ws.Range("F13").Value2 = $periods.Count                 # Cant. Periodos
$ws.Range("E11").Value2 = 3265556                         # VALOR MORA total
